$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.233.38'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.79%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.453.72'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.15%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '569.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.95%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.99'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.84%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +0.98%  '
$ws.Range('E9').Value = '  +12.43%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '2.452.50'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.10%  '
$ws.Range('E11').Value = '  -1.68%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.337'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.77%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.70'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.70%  '
$ws.Range('E14').Value = '  +8.08%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '70.114.84'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.73%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.904.74'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.20'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.47%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.447.56'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.90'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.56%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.16'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.01%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '341.96'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.56%  '
$ws.Range('E22').Value = '  +3.40%  '
$ws.Range('E23').Value = '  +8.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.53'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.59%  '
$ws.Range('E26').Value = '  +6.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.56'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.97%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.579.80'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.991'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.76%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0859'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.71%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.40'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.64%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '465.05'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +11.07%  '
$ws.Range('E33').Value = '  +10.74%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.63'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '159.38'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.95%  '
$ws.Range('E37').Value = '  +7.71%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.12'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.24'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.55%  '
$ws.Range('E42').Value = '  +5.17%  '
$ws.Range('E43').Value = '  +4.41%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '38.12'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.87%  '
$ws.Range('E45').Value = '  +3.28%  '
$ws.Range('E46').Value = '  +6.65%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '134.44'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.55%  '
$ws.Range('E48').Value = '  +2.28%  '
$ws.Range('E49').Value = '  +2.77%  '
$ws.Range('E50').Value = '  +3.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.565'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.08%  '
